# Applies the "LinuxForHealth" re-brand edit described in the commit:
#  - Metadata sheet: URL / Version / Date / Publisher text updates.
#  - Elements sheet: the shared-string table entries for the
#    "Extension.id" and "Extension.extension" element rows were
#    reordered ahead of the inherited "ele-1/ext-1" constraint text,
#    which is why that constraint text now surfaces on the
#    Extension.value[x] row instead of the root Extension row, and a
#    cluster of cells on rows 3-6 pick up shifted values. Expressed
#    here purely as literal cell-value writes (no rows/cols altered).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Metadata"
# ---------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/stated-reason"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---------------------------------------------------------------
# Sheet "Elements"
# ---------------------------------------------------------------
$el = $wb.Worksheets.Item("Elements")

$el.Range("AI2").Value = ""
$el.Range("A3").Value = "Extension.id"
$el.Range("AE3").Value = "Element.id"
$el.Range("AG3").Value = "1"
$el.Range("AJ3").Value = "n/a"
$el.Range("F3").Value = "1"
$el.Range("J3").Value = "string`n"
$el.Range("K3").Value = "Unique id for inter-element referencing"
$el.Range("L3").Value = "Unique id for the element within a resource (for internal references). This may be any string value that does not contain spaces."
$el.Range("A4").Value = "Extension.extension"
$el.Range("AA4").Value = "value:url}`n"
$el.Range("AB4").Value = "Extensions are always sliced by (at least) url"
$el.Range("AD4").Value = "open"
$el.Range("AE4").Value = "Element.extension"
$el.Range("AI4").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
$el.Range("J4").Value = "Extension`n"
$el.Range("L4").Value = "An Extension"
$el.Range("AF5").Value = "1"
$el.Range("AG5").Value = "1"
$el.Range("E5").Value = "1"
$el.Range("F5").Value = "1"
$el.Range("AG6").Value = "1"
$el.Range("F6").Value = "1"

